$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataSetSections")
$ws.Range("A2").Value2 = "SETVIA_VALUE2"
